$d = $word.ActiveDocument

$dquote = [char]0x201C
$rquote = [char]0x201D

# --- Step 1: split the sub-bullet paragraph in two, right after the
# existing sentence, inserting the new answer text "Yes, yes I have. " as
# its own paragraph. Doing the split first (while the paragraph is still
# plain/unformatted) keeps the new paragraph free of inherited direct
# formatting. ---

$searchText  = "Triangle Effect." + $rquote
$replaceText = "Triangle Effect." + $rquote + "`rYes, yes I have. "
$d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replaceText, 2) | Out-Null

# --- Step 2: the original sentence (now alone in its paragraph, still as 3
# runs with no formatting) becomes a single strikethrough run with the same
# combined text. Rebuilding the paragraph's text collapses it to one run;
# applying strikethrough afterwards formats the whole paragraph. ---

$p = $d.Paragraphs.Item(8)
$fullText = "I verified that the Euler Cromer Method and the Leapfrog Algorithm" `
    + " both converge as approcimately O(h) I have no clue why my leapfrog" `
    + " algorithm is acting almost just like the Euler Cromer. Instead I have" `
    + " discovered the " + $dquote + "Triangle Effect." + $rquote

$rng = $d.Range($p.Range.Start, $p.Range.End - 1)
$rng.Delete()
$rng.InsertAfter($fullText)

$p = $d.Paragraphs.Item(8)
$p.Range.Font.StrikeThrough = 1

# --- Step 3: the new answer paragraph loses its inherited numbering (it
# becomes a plain ListParagraph, not a sub-bullet) and has strike/dstrike
# explicitly cleared. ---

$p2 = $d.Paragraphs.Item(9)
$p2.Style = "ListParagraph"
$p2.Range.Font.StrikeThrough = 0
$p2.Range.Font.DoubleStrikeThrough = 0

Write-Output "done"
